# Apply the WM2018 knock-out stage updates (Round of 16 through the 3rd-place game)
# to the "Matches" sheet, plus the cursor/selection state change recorded in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Matches")

# --- Seed the new shared strings in the exact order they were first typed, using a
#     scratch column far outside the used range so the shared-string table ends up
#     in the same order as the source workbook. The scratch cells are cleared again
#     at the end, once every string is referenced by the real data cells below. ---
$ws.Range("ZZ1").Value = "Sieger C"
$ws.Range("ZZ2").Value = "Zweiter D"
$ws.Range("ZZ3").Value = "Sieger A"
$ws.Range("ZZ4").Value = "Zweiter B"
$ws.Range("ZZ5").Value = "Sieger B"
$ws.Range("ZZ6").Value = "Zweiter A"
$ws.Range("ZZ7").Value = "Sieger D"
$ws.Range("ZZ8").Value = "Zweiter C"
$ws.Range("ZZ9").Value = "Sieger E"
$ws.Range("ZZ10").Value = "Zweiter F"
$ws.Range("ZZ11").Value = "Sieger G"
$ws.Range("ZZ12").Value = "Zweiter H"
$ws.Range("ZZ13").Value = "Sieger F"
$ws.Range("ZZ14").Value = "Sieger H"
$ws.Range("ZZ15").Value = "Zweiter E"
$ws.Range("ZZ16").Value = "Zweiter G"
$ws.Range("ZZ17").Value = "Sieger A1"
$ws.Range("ZZ18").Value = "Zweiter A2"
$ws.Range("ZZ19").Value = "Sieger A3"
$ws.Range("ZZ20").Value = "Zweiter A4"
$ws.Range("ZZ21").Value = "Sieger A5"
$ws.Range("ZZ22").Value = "Zweiter A6"
$ws.Range("ZZ23").Value = "Sieger A7"
$ws.Range("ZZ24").Value = "Zweiter A8"
$ws.Range("ZZ25").Value = "Sieger HF1"
$ws.Range("ZZ26").Value = "Sieger HF2"
$ws.Range("ZZ27").Value = "Verlierer HF1"
$ws.Range("ZZ28").Value = "Verlierer HF2"
$ws.Range("ZZ29").Value = "Sieger VF1"
$ws.Range("ZZ30").Value = "Sieger VF2"
$ws.Range("ZZ31").Value = "Sieger VF4"
$ws.Range("ZZ32").Value = "Sieger VF3"

# --- Round of 16 / Quarter-final / Semi-final / Final / 3rd-place fixtures ---
$ws.Range("A50").Value = "Sieger C"
$ws.Range("B50").Value = "Zweiter D"
$ws.Range("C50").Value = "ROUND_OF_SIXTEEN"
$ws.Range("D50").Value = 43281.666666666664
$ws.Range("E50").Value = "Kasan"

$ws.Range("A51").Value = "Sieger A"
$ws.Range("B51").Value = "Zweiter B"
$ws.Range("C51").Value = "ROUND_OF_SIXTEEN"
$ws.Range("D51").Value = 43281.833333333336
$ws.Range("E51").Value = "Sotschi"

$ws.Range("A52").Value = "Sieger B"
$ws.Range("B52").Value = "Zweiter A"
$ws.Range("C52").Value = "ROUND_OF_SIXTEEN"
$ws.Range("D52").Value = 43282.75
$ws.Range("E52").Value = "Moskau"

$ws.Range("A53").Value = "Sieger D"
$ws.Range("B53").Value = "Zweiter C"
$ws.Range("C53").Value = "ROUND_OF_SIXTEEN"
$ws.Range("D53").Value = 43282.833333333336
$ws.Range("E53").Value = "Nischni Nowgorod"

$ws.Range("A54").Value = "Sieger E"
$ws.Range("B54").Value = "Zweiter F"
$ws.Range("C54").Value = "ROUND_OF_SIXTEEN"
$ws.Range("D54").Value = 43283.666666666664
$ws.Range("E54").Value = "Samara"

$ws.Range("A55").Value = "Sieger G"
$ws.Range("B55").Value = "Zweiter H"
$ws.Range("C55").Value = "ROUND_OF_SIXTEEN"
$ws.Range("D55").Value = 43283.833333333336
$ws.Range("E55").Value = "Rostow am Don"

$ws.Range("A56").Value = "Sieger F"
$ws.Range("B56").Value = "Zweiter E"
$ws.Range("C56").Value = "ROUND_OF_SIXTEEN"
$ws.Range("D56").Value = 43284.666666666664
$ws.Range("E56").Value = "Sankt Petersburg"

$ws.Range("A57").Value = "Sieger H"
$ws.Range("B57").Value = "Zweiter G"
$ws.Range("C57").Value = "ROUND_OF_SIXTEEN"
$ws.Range("D57").Value = 43284.833333333336
$ws.Range("E57").Value = "Moskau"

$ws.Range("A58").Value = "Sieger A1"
$ws.Range("B58").Value = "Zweiter A2"
$ws.Range("C58").Value = "QUARTER_FINAL"
$ws.Range("D58").Value = 43287.666666666664
$ws.Range("E58").Value = "Nischni Nowgorod"

$ws.Range("A59").Value = "Sieger A3"
$ws.Range("B59").Value = "Zweiter A4"
$ws.Range("C59").Value = "QUARTER_FINAL"
$ws.Range("D59").Value = 43288.833333333336
$ws.Range("E59").Value = "Sotschi"

$ws.Range("A60").Value = "Sieger A5"
$ws.Range("B60").Value = "Zweiter A6"
$ws.Range("C60").Value = "QUARTER_FINAL"
$ws.Range("D60").Value = 43287.833333333336
$ws.Range("E60").Value = "Kasan"

$ws.Range("A61").Value = "Sieger A7"
$ws.Range("B61").Value = "Zweiter A8"
$ws.Range("C61").Value = "QUARTER_FINAL"
$ws.Range("D61").Value = 43288.666666666664
$ws.Range("E61").Value = "Samara"

$ws.Range("A62").Value = "Sieger VF2"
$ws.Range("B62").Value = "Sieger VF1"
$ws.Range("C62").Value = "SEMI_FINAL"
$ws.Range("D62").Value = 43291.833333333336
$ws.Range("E62").Value = "Sankt Petersburg"

$ws.Range("A63").Value = "Sieger VF4"
$ws.Range("B63").Value = "Sieger VF3"
$ws.Range("C63").Value = "SEMI_FINAL"
$ws.Range("D63").Value = 43292.833333333336
$ws.Range("E63").Value = "Moskau"

$ws.Range("A64").Value = "Sieger HF1"
$ws.Range("B64").Value = "Sieger HF2"
$ws.Range("C64").Value = "FINAL"
$ws.Range("D64").Value = 43296.708333333336
$ws.Range("E64").Value = "Moskau"

$ws.Range("A65").Value = "Verlierer HF1"
$ws.Range("B65").Value = "Verlierer HF2"
$ws.Range("C65").Value = "GAME_FOR_THIRD"
$ws.Range("D65").Value = 43295.666666666664
$ws.Range("E65").Value = "Sankt Petersburg"

# --- Clear the scratch cells now that every new string has a real reference ---
$ws.Range("ZZ1").Value = $null
$ws.Range("ZZ2").Value = $null
$ws.Range("ZZ3").Value = $null
$ws.Range("ZZ4").Value = $null
$ws.Range("ZZ5").Value = $null
$ws.Range("ZZ6").Value = $null
$ws.Range("ZZ7").Value = $null
$ws.Range("ZZ8").Value = $null
$ws.Range("ZZ9").Value = $null
$ws.Range("ZZ10").Value = $null
$ws.Range("ZZ11").Value = $null
$ws.Range("ZZ12").Value = $null
$ws.Range("ZZ13").Value = $null
$ws.Range("ZZ14").Value = $null
$ws.Range("ZZ15").Value = $null
$ws.Range("ZZ16").Value = $null
$ws.Range("ZZ17").Value = $null
$ws.Range("ZZ18").Value = $null
$ws.Range("ZZ19").Value = $null
$ws.Range("ZZ20").Value = $null
$ws.Range("ZZ21").Value = $null
$ws.Range("ZZ22").Value = $null
$ws.Range("ZZ23").Value = $null
$ws.Range("ZZ24").Value = $null
$ws.Range("ZZ25").Value = $null
$ws.Range("ZZ26").Value = $null
$ws.Range("ZZ27").Value = $null
$ws.Range("ZZ28").Value = $null
$ws.Range("ZZ29").Value = $null
$ws.Range("ZZ30").Value = $null
$ws.Range("ZZ31").Value = $null
$ws.Range("ZZ32").Value = $null

# --- Restore the view state recorded for this sheet: scrolled down with B63 selected ---
$ws.Activate()
$ws.Range("B63").Select()

